$d = $word.ActiveDocument

# 1. Apply the "No Spacing" paragraph style to every paragraph in the body.
#    (This both stamps <w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr> on each
#    paragraph and mints the NoSpacing style definition in styles.xml.)
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $d.Paragraphs($i).Style = "No Spacing"
}

# 2. Tidy up the auto-generated "No Spacing" style definition so it matches
#    the canonical Word "No Spacing" quick style (uiPriority 1, no
#    w:basedOn, single line-spacing/no space-after pPr).
$styles = $d.Styles
$ns = $styles.Item("No Spacing")
$ns.Priority = 1
$ns.ParagraphFormat.LineSpacingRule = 0   # wdLineSpaceSingle
$ns.ParagraphFormat.SpaceAfter = 0

# 3. Move the "_GoBack" bookmark from the end of the "...week" paragraph to
#    the end of the final "...for Edmond" paragraph, so it now brackets the
#    very end of the document instead of sitting between the two paragraphs.
$last = $d.Paragraphs($count).Range
$last.End = $last.End - 1   # exclude the paragraph mark
$last.Collapse(0)           # wdCollapseEnd
$last.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $last)
$bk = $d.Bookmarks.Item("_GoBack")
$bk.Range.Delete()
